$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host "D width:" $ws.Columns.Item(4).ColumnWidth
Write-Host "E width:" $ws.Columns.Item(5).ColumnWidth
Write-Host "F width:" $ws.Columns.Item(6).ColumnWidth
